# Apply targeted updates to column A (rowid) values in Sheet1
# These correspond to shared-string value corrections from the source commit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    99 = "obs20790"
    100 = "obs20972"
    101 = "obs21165"
    205 = "obs20599"
    206 = "obs20792"
    207 = "obs20974"
    208 = "obs21167"
    224 = "obs22146"
    225 = "obs22352"
    226 = "obs22574"
    227 = "obs22810"
    228 = "obs23057"
    229 = "obs23323"
    230 = "obs23605"
    231 = "obs23914"
    232 = "obs24238"
    233 = "obs24565"
    234 = "obs24913"
    235 = "obs25275"
    236 = "obs25652"
    237 = "obs26061"
    238 = "obs22147"
    239 = "obs22353"
    240 = "obs22575"
    241 = "obs22811"
    242 = "obs23105"
    243 = "obs23324"
    244 = "obs23606"
    245 = "obs23915"
    246 = "obs24239"
    247 = "obs24566"
    248 = "obs24914"
    249 = "obs25276"
    250 = "obs25653"
    251 = "obs26062"
    252 = "obs22354"
    253 = "obs22576"
    254 = "obs22812"
    255 = "obs23058"
    256 = "obs23325"
    257 = "obs23607"
    258 = "obs23916"
    259 = "obs24240"
    260 = "obs24567"
    261 = "obs24915"
    262 = "obs25277"
    263 = "obs25654"
    264 = "obs26063"
    265 = "obs22148"
    266 = "obs22355"
    267 = "obs22577"
    268 = "obs22813"
    269 = "obs23059"
    270 = "obs23326"
    271 = "obs23608"
    272 = "obs23917"
    273 = "obs24241"
    274 = "obs24568"
    275 = "obs24916"
    276 = "obs25278"
    277 = "obs25655"
    278 = "obs26064"
    279 = "obs22153"
    280 = "obs22360"
    281 = "obs22582"
    282 = "obs22816"
    283 = "obs23062"
    284 = "obs23329"
    285 = "obs23611"
    286 = "obs23919"
    287 = "obs24297"
    288 = "obs24570"
    289 = "obs24918"
    290 = "obs25280"
    291 = "obs25657"
    292 = "obs26066"
    293 = "obs22853"
    294 = "obs23107"
    295 = "obs23371"
    296 = "obs23612"
    297 = "obs23920"
    298 = "obs24299"
    299 = "obs24571"
    300 = "obs24919"
    301 = "obs25281"
    302 = "obs25658"
    303 = "obs26070"
    304 = "obs22855"
    305 = "obs23109"
    306 = "obs23373"
    307 = "obs23613"
    308 = "obs23921"
    309 = "obs24301"
    310 = "obs24572"
    311 = "obs24920"
    312 = "obs25282"
    313 = "obs25659"
    314 = "obs26067"
    315 = "obs23661"
    316 = "obs23979"
    317 = "obs24302"
    318 = "obs24629"
    319 = "obs24981"
    320 = "obs25343"
    321 = "obs25726"
    322 = "obs26160"
    323 = "obs23662"
    324 = "obs23980"
    325 = "obs24303"
    326 = "obs24630"
    327 = "obs24982"
    328 = "obs25344"
    329 = "obs25727"
    330 = "obs26161"
    331 = "obs23663"
    332 = "obs23981"
    333 = "obs24304"
    334 = "obs24631"
    335 = "obs24983"
    336 = "obs25345"
    337 = "obs25728"
    338 = "obs26162"
    339 = "obs25680"
    340 = "obs26100"
    341 = "obs25681"
    342 = "obs26101"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}

